$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Header/Footer "fixed" date placeholders (slide master + every slide
#    layout) were re-stamped from 29/01/2023 to 31/01/2023.
# ---------------------------------------------------------------------------
function Update-DateShape($shapes) {
    for ($j = 1; $j -le $shapes.Count; $j++) {
        $shp = $shapes.Item($j)
        if ($shp.Name -like "Date Placeholder*") {
            $shp.TextFrame.TextRange.Text = "31/01/2023"
        }
    }
}

$design    = $p.Designs.Item(1)
$masterObj = $design.SlideMaster

Update-DateShape $masterObj.Shapes

$layouts = $masterObj.CustomLayouts
for ($i = 1; $i -le $layouts.Count; $i++) {
    Update-DateShape $layouts.Item($i).Shapes
}

# ---------------------------------------------------------------------------
# 2) Slide 4 ("Tools Used") content placeholder gained three extra bullet
#    lines describing the tech stack used for the routes.
# ---------------------------------------------------------------------------
$toolsSlide = $p.Slides.Item(4)
$content    = $toolsSlide.Shapes.Item(2)
$content.TextFrame.TextRange.Text = "Slack`rZoom`rJIRA`rGITHUB`rPython`rFlask`rMongo(will be updated to SQL)"
